$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 32.8
$ws.Range("I8").Value = 16
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 48
$ws.Range("L8").Value = 300
$ws.Range("M8").Value = 91
$ws.Range("N8").Value = -578

$ws.Range("H28").Value = 9453.416999999999
$ws.Range("I28").Value = 14650.286
$ws.Range("J28").Value = 2177.8
$ws.Range("K28").Value = 14650.286
$ws.Range("L28").Value = 2177.8
$ws.Range("M28").Value = -14165.286
$ws.Range("N28").Value = -3147.8

$ws.Range("H95").Value = 37499.5
$ws.Range("J95").Value = 37499.5
$ws.Range("L95").Value = 37499.5
$ws.Range("N95").Value = -42991.5

$ws.Range("H100").Value = 3666.6667
$ws.Range("I100").Value = 3666.6667
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3666.6667
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3125.6667

$ws.Range("H125").Value = 4975.636
$ws.Range("I125").Value = 2354.125
$ws.Range("J125").Value = 11966.333
$ws.Range("K125").Value = 21187.125
$ws.Range("L125").Value = 107696.997
$ws.Range("M125").Value = -18727.125
$ws.Range("N125").Value = -112616.997

$ws.Range("H131").Value = 299.5
$ws.Range("I131").Value = 299.5
$ws.Range("K131").Value = 898.5
$ws.Range("M131").Value = 4141.5

$ws.Range("H132").Value = 5125.6665
$ws.Range("I132").Value = 5364.364
$ws.Range("K132").Value = 16093.092
$ws.Range("M132").Value = -13563.092

$ws.Range("H138").Value = 2799.25
$ws.Range("I138").Value = 2197
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 6591
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -1451
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 35863.75
$ws.Range("J24").Value = 35863.75
$ws.Range("L24").Value = 35863.75
$ws.Range("N24").Value = -36611.75

$ws.Range("H92").Value = 36633.332
$ws.Range("J92").Value = 36633.332
$ws.Range("L92").Value = 36633.332
$ws.Range("N92").Value = -41625.332

$ws.Range("H100").Value = 35863.75
$ws.Range("J100").Value = 35863.75
$ws.Range("L100").Value = 35863.75
$ws.Range("N100").Value = -38027.75

$ws.Range("H110").Value = 1677.5625
$ws.Range("I110").Value = 1503.9231
$ws.Range("K110").Value = 1503.9231
$ws.Range("M110").Value = 541.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5759
$ws.Range("I20").Value = 5365.1
$ws.Range("K20").Value = 5365.1
$ws.Range("M20").Value = -5118.1

$ws.Range("H88").Value = 15389.4
$ws.Range("J88").Value = 15389.4
$ws.Range("L88").Value = 15389.4
$ws.Range("N88").Value = -16201.4

$ws.Range("H91").Value = 15389.4
$ws.Range("J91").Value = 15389.4
$ws.Range("L91").Value = 15389.4
$ws.Range("N91").Value = -18197.4

$ws.Range("H99").Value = 1042.5555
$ws.Range("I99").Value = 1042.5555
$ws.Range("K99").Value = 1042.5555
$ws.Range("M99").Value = 455.4445000000001

$ws.Range("H105").Value = 1816.25
$ws.Range("I105").Value = 1217.8889
$ws.Range("J105").Value = 2585.5715
$ws.Range("K105").Value = 1217.8889
$ws.Range("L105").Value = 2585.5715
$ws.Range("M105").Value = 529.1111000000001
$ws.Range("N105").Value = -6079.5715

$ws.Range("H107").Value = 934.75
$ws.Range("I107").Value = 782.5714
$ws.Range("K107").Value = 782.5714
$ws.Range("M107").Value = 1137.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1137.7142
$ws.Range("I16").Value = 1137.7142
$ws.Range("K16").Value = 1137.7142
$ws.Range("M16").Value = -850.7141999999999

$ws.Range("H99").Value = 5504.273
$ws.Range("I99").Value = 5506
$ws.Range("J99").Value = 5499.6665
$ws.Range("K99").Value = 5506
$ws.Range("L99").Value = 5499.6665
$ws.Range("M99").Value = -4008
$ws.Range("N99").Value = -8495.666499999999

$ws.Range("H103").Value = 3000
$ws.Range("I103").Value = 3000
$ws.Range("K103").Value = 3000
$ws.Range("M103").Value = -1828

$ws.Range("H113").Value = 1137.7142
$ws.Range("I113").Value = 1137.7142
$ws.Range("K113").Value = 1137.7142
$ws.Range("M113").Value = 1032.2858

$ws.Range("H126").Value = 5504.273
$ws.Range("I126").Value = 5506
$ws.Range("J126").Value = 5499.6665
$ws.Range("K126").Value = 16518
$ws.Range("L126").Value = 16498.9995
$ws.Range("M126").Value = -14048
$ws.Range("N126").Value = -21438.9995

$ws.Range("H134").Value = 1527.4762
$ws.Range("I134").Value = 1579.3334
$ws.Range("K134").Value = 4738.0002
$ws.Range("M134").Value = -2203.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 621.75
$ws.Range("I11").Value = 755.75
$ws.Range("K11").Value = 2267.25
$ws.Range("M11").Value = -2127.25

$ws.Range("H98").Value = 198.5
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H137").Value = 1030
$ws.Range("I137").Value = 1030
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3090
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = 2010

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6147
$ws.Range("I70").Value = 5997
$ws.Range("J70").Value = 6297
$ws.Range("K70").Value = 5997
$ws.Range("L70").Value = 6297
$ws.Range("M70").Value = -5727
$ws.Range("N70").Value = -6837

$ws.Range("H73").Value = 6147
$ws.Range("I73").Value = 5997
$ws.Range("J73").Value = 6297
$ws.Range("K73").Value = 5997
$ws.Range("L73").Value = 6297
$ws.Range("M73").Value = -5061
$ws.Range("N73").Value = -8169

$ws.Range("H113").Value = 4170
$ws.Range("I113").Value = 4170
$ws.Range("K113").Value = 4170
$ws.Range("M113").Value = -2000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1710.5555
$ws.Range("J22").Value = 2950
$ws.Range("L22").Value = 2950
$ws.Range("N22").Value = -3540

$ws.Range("H27").Value = 1710.5555
$ws.Range("J27").Value = 2950
$ws.Range("L27").Value = 2950
$ws.Range("N27").Value = -3164

$ws.Range("H46").Value = 4454.5454
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H55").Value = 305.6842
$ws.Range("I55").Value = 259.63635
$ws.Range("K55").Value = 259.63635
$ws.Range("M55").Value = -86.63634999999999

$ws.Range("H68").Value = 24320
$ws.Range("I68").Value = 2900
$ws.Range("K68").Value = 2900
$ws.Range("M68").Value = -2151

$ws.Range("H71").Value = 24320
$ws.Range("I71").Value = 2900
$ws.Range("K71").Value = 14500
$ws.Range("M71").Value = -10756

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0

$ws.Range("H93").Value = 874.4
$ws.Range("I93").Value = 874.4
$ws.Range("K93").Value = 874.4
$ws.Range("M93").Value = 373.6

$ws.Range("H132").Value = 6026.125
$ws.Range("I132").Value = 5740.8
$ws.Range("J132").Value = 6155.8184
$ws.Range("K132").Value = 17222.4
$ws.Range("L132").Value = 18467.4552
$ws.Range("M132").Value = -14692.4
$ws.Range("N132").Value = -23527.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4049
$ws.Range("I62").Value = 4064.8333
$ws.Range("J62").Value = 4001.5
$ws.Range("K62").Value = 4064.8333
$ws.Range("L62").Value = 4001.5
$ws.Range("M62").Value = -3440.8333
$ws.Range("N62").Value = -5249.5

$ws.Range("H65").Value = 4049
$ws.Range("I65").Value = 4064.8333
$ws.Range("J65").Value = 4001.5
$ws.Range("K65").Value = 20324.1665
$ws.Range("L65").Value = 20007.5
$ws.Range("M65").Value = -17204.1665
$ws.Range("N65").Value = -26247.5

$ws.Range("H95").Value = 27200
$ws.Range("J95").Value = 27200
$ws.Range("L95").Value = 27200
$ws.Range("N95").Value = -32692

$ws.Range("H113").Value = 7949.7144
$ws.Range("J113").Value = 838.3333
$ws.Range("L113").Value = 2514.9999
$ws.Range("N113").Value = -6854.9999

$ws.Range("H132").Value = 1767.8422
$ws.Range("I132").Value = 1611.3529
$ws.Range("K132").Value = 4834.0587
$ws.Range("M132").Value = -2304.0587
